# Add three new Amazon order rows (46, 47, 48) to the bottom of the sheet,
# mirroring the existing "Raw Sales Data" layout (columns B, I, S-X, AC, AE, AF
# are intentionally left blank, same as every prior row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 46 ----------------------------------------------------------
$ws.Cells.Item(46, 1).Value  = "026-8963265-6046736"                 # A amazon-order-id
$ws.Cells.Item(46, 3).Value  = "2025-09-04T18:44:57+00:00"           # C purchase-date
$ws.Cells.Item(46, 4).Value  = "2025-09-04T19:14:34+00:00"           # D last-updated-date
$ws.Cells.Item(46, 5).Value  = "Pending"                             # E order-status
$ws.Cells.Item(46, 6).Value  = "Merchant"                            # F fulfillment-channel
$ws.Cells.Item(46, 7).Value  = "Amazon.co.uk"                        # G sales-channel
$ws.Cells.Item(46, 8).Value  = "WebsiteOrderChannel"                 # H order-channel
$ws.Cells.Item(46, 10).Value = "Standard"                            # J ship-service-level
$ws.Cells.Item(46, 11).Value = $ws.Cells.Item(7, 11).Value()         # K product-name (reused text)
$ws.Cells.Item(46, 12).Value = "Coffee Table White 250601000006"     # L sku
$ws.Cells.Item(46, 13).Value = "B0FLQN7PRC"                          # M asin
$ws.Cells.Item(46, 14).Value = "Unshipped"                           # N item-status
$ws.Cells.Item(46, 15).Value = 1                                     # O quantity
$ws.Cells.Item(46, 16).Value = "GBP"                                 # P currency
$ws.Cells.Item(46, 17).Value = 34.99                                 # Q item-price
$ws.Cells.Item(46, 18).Value = 5.83                                  # R item-tax
$ws.Cells.Item(46, 25).Value = "Coppull Chorley"                     # Y ship-city
$ws.Cells.Item(46, 26).Value = "Lancs"                                # Z ship-state
$ws.Cells.Item(46, 27).Value = "PR7 4PA"                             # AA ship-postal-code
$ws.Cells.Item(46, 28).Value = "GB"                                  # AB ship-country
$ws.Cells.Item(46, 30).Value = $false                                # AD is-business-order
$ws.Cells.Item(46, 33).Value = $false                                # AG is-iba

# ---- Row 47 ----------------------------------------------------------
$ws.Cells.Item(47, 1).Value  = "203-8858258-7768358"                 # A amazon-order-id
$ws.Cells.Item(47, 3).Value  = "2025-09-04T16:36:07+00:00"           # C purchase-date
$ws.Cells.Item(47, 4).Value  = "2025-09-04T23:26:06+00:00"           # D last-updated-date
$ws.Cells.Item(47, 5).Value  = "Pending"                             # E order-status
$ws.Cells.Item(47, 6).Value  = "Merchant"                            # F fulfillment-channel
$ws.Cells.Item(47, 7).Value  = "Amazon.co.uk"                        # G sales-channel
$ws.Cells.Item(47, 8).Value  = "WebsiteOrderChannel"                 # H order-channel
$ws.Cells.Item(47, 10).Value = "Standard"                            # J ship-service-level
$ws.Cells.Item(47, 11).Value = $ws.Cells.Item(4, 11).Value()         # K product-name (reused text)
$ws.Cells.Item(47, 12).Value = "Executive Black Chair 250601000002"  # L sku
$ws.Cells.Item(47, 13).Value = "B0FLQMG769"                          # M asin
$ws.Cells.Item(47, 14).Value = "Unshipped"                           # N item-status
$ws.Cells.Item(47, 15).Value = 1                                     # O quantity
$ws.Cells.Item(47, 16).Value = "GBP"                                 # P currency
$ws.Cells.Item(47, 17).Value = 99                                    # Q item-price
$ws.Cells.Item(47, 18).Value = 16.5                                  # R item-tax
$ws.Cells.Item(47, 25).Value = "Nottingham"                          # Y ship-city
$ws.Cells.Item(47, 26).Value = "Nottinghamshire"                     # Z ship-state
$ws.Cells.Item(47, 27).Value = "NG2 7GL"                             # AA ship-postal-code
$ws.Cells.Item(47, 28).Value = "GB"                                  # AB ship-country
$ws.Cells.Item(47, 30).Value = $false                                # AD is-business-order
$ws.Cells.Item(47, 33).Value = $false                                # AG is-iba

# ---- Row 48 ----------------------------------------------------------
$ws.Cells.Item(48, 1).Value  = "202-4176722-0809902"                 # A amazon-order-id
$ws.Cells.Item(48, 3).Value  = "2025-09-04T06:20:25+00:00"           # C purchase-date
$ws.Cells.Item(48, 4).Value  = "2025-09-04T08:08:32+00:00"           # D last-updated-date
$ws.Cells.Item(48, 5).Value  = "Shipped"                             # E order-status
$ws.Cells.Item(48, 6).Value  = "Merchant"                            # F fulfillment-channel
$ws.Cells.Item(48, 7).Value  = "Amazon.co.uk"                        # G sales-channel
$ws.Cells.Item(48, 8).Value  = "WebsiteOrderChannel"                 # H order-channel
$ws.Cells.Item(48, 10).Value = "Standard"                            # J ship-service-level
$ws.Cells.Item(48, 11).Value = $ws.Cells.Item(2, 11).Value()         # K product-name (reused text)
$ws.Cells.Item(48, 12).Value = "Electric Drill 2506010000001"        # L sku
$ws.Cells.Item(48, 13).Value = "B0FLQKQLX5"                          # M asin
$ws.Cells.Item(48, 14).Value = "Shipped"                             # N item-status
$ws.Cells.Item(48, 15).Value = 1                                     # O quantity
$ws.Cells.Item(48, 16).Value = "GBP"                                 # P currency
$ws.Cells.Item(48, 17).Value = 23.99                                 # Q item-price
$ws.Cells.Item(48, 18).Value = 4                                     # R item-tax
$ws.Cells.Item(48, 25).Value = "LIVERPOOL"                           # Y ship-city
$ws.Cells.Item(48, 27).Value = "L10 0AG"                             # AA ship-postal-code
$ws.Cells.Item(48, 28).Value = "GB"                                  # AB ship-country
$ws.Cells.Item(48, 30).Value = $false                                # AD is-business-order
$ws.Cells.Item(48, 33).Value = $false                                # AG is-iba

# ---- View state to mirror the saved workbook --------------------------
# Scroll so row 37 is the top visible row (was row 22 before the edit),
# then leave the selection on D52 like the source file.
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D52").Select()
